$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23:86 down to 24:87
$ws.Rows(23).Insert()

# Populate the newly inserted row 23 with the new weekly price-report entry
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value = 44544
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100112008
$ws.Range("G23").Value = "Coliflor"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 1200
$ws.Range("K23").Value = 400
$ws.Range("L23").Value = 500
$ws.Range("M23").Value = 450
$ws.Range("N23").Value = "`$/unidad"
$ws.Range("O23").Value = "Región de Arica y Parinacota"
$ws.Range("P23").Value = 450
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
